# modified TB mort and contact rates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_matched_parameters")

# --- Contact rate (beta) inputs, row 2-3 ---
# J/K columns hold =I*0.75 / =I*1.25 formulas and recompute automatically.
$ws.Range("I2").Value = 12    # beta_1 : 11 -> 12
$ws.Range("I3").Value = 12    # beta_2 : 11 -> 12

# --- TB mortality / risk inputs, rows 38-39, 44-45 ---
$ws.Range("I38").Value = 28   # risk.other_3 : 30   -> 28
$ws.Range("I39").Value = 1.34 # risk.other_4 : 1.35 -> 1.34
$ws.Range("I44").Value = 26   # risk.TB_2    : 27   -> 26
$ws.Range("I45").Value = 52   # risk.TB_3    : 51   -> 52

# Row 37 was manually resized (wrapped-text autofit row shrank from 34pt to 17pt)
$ws.Rows.Item(37).RowHeight = 17

# Restore the author's last cursor position / selection on the sheet
$ws.Activate()
$ws.Range("F9").Select()

# Reposition the application window on screen to match the saved workbook view
$win = $excel.ActiveWindow
$win.Left = 2720
$win.Top = 1180

$wb.Save()
